$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.660.94"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "3.504.07"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "3.503.98"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "4.103.19"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "3.508.20"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "64.746.54"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "3.645.95"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -3.80%  "
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.14%  "
$ws.Range("D33").Value = "3.506.91"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +3.60%  "
$ws.Range("E38").Value = "  +4.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "171.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0815"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.814"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").Value = "2.475.50"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("E51").Value = "  +4.54%  "
